$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 50003.74
$ws.Range("I28").Value = 1339.25
$ws.Range("J28").Value = 133428.58
$ws.Range("K28").Value = 1339.25
$ws.Range("L28").Value = 133428.58
$ws.Range("M28").Value = -854.25
$ws.Range("N28").Value = -134398.58
$ws.Range("H33").Value = 1649.7693
$ws.Range("I33").Value = 1649.7693
$ws.Range("K33").Value = 1649.7693
$ws.Range("M33").Value = -1420.7693
$ws.Range("H40").Value = 8401.883
$ws.Range("I40").Value = 2318.7778
$ws.Range("K40").Value = 2318.7778
$ws.Range("M40").Value = -2143.7778
$ws.Range("H88").Value = 1502439.8
$ws.Range("J88").Value = 3000
$ws.Range("L88").Value = 3000
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 1502439.8
$ws.Range("J91").Value = 3000
$ws.Range("L91").Value = 3000
$ws.Range("N91").Value = -5808
$ws.Range("H125").Value = 4315.75
$ws.Range("I125").Value = 3999.6667
$ws.Range("J125").Value = 4505.4
$ws.Range("K125").Value = 35997.0003
$ws.Range("L125").Value = 40548.6
$ws.Range("M125").Value = -33537.0003
$ws.Range("N125").Value = -45468.6
$ws.Range("H129").Value = 1950.8695
$ws.Range("I129").Value = 1477
$ws.Range("J129").Value = 3656.8
$ws.Range("K129").Value = 4431
$ws.Range("L129").Value = 10970.4
$ws.Range("M129").Value = 569
$ws.Range("N129").Value = -20970.4
$ws.Range("H132").Value = 1323.1719
$ws.Range("I132").Value = 1196.2106
$ws.Range("J132").Value = 2357
$ws.Range("K132").Value = 3588.6318
$ws.Range("L132").Value = 7071
$ws.Range("M132").Value = -1058.6318
$ws.Range("N132").Value = -12131
$ws.Range("H134").Value = 91897.375
$ws.Range("J134").Value = 91897.375
$ws.Range("L134").Value = 91897.375
$ws.Range("N134").Value = -102037.375
$ws.Range("H137").Value = 324710.06
$ws.Range("I137").Value = 2455.8
$ws.Range("J137").Value = 1452600
$ws.Range("K137").Value = 7367.400000000001
$ws.Range("L137").Value = 4357800
$ws.Range("M137").Value = -4817.400000000001
$ws.Range("N137").Value = -4362900
$ws.Range("H138").Value = 2334.2683
$ws.Range("I138").Value = 1562.5714
$ws.Range("J138").Value = 2734.4075
$ws.Range("K138").Value = 4687.7142
$ws.Range("L138").Value = 8203.2225
$ws.Range("M138").Value = 452.2857999999997
$ws.Range("N138").Value = -18483.2225
$ws.Range("H140").Value = 70487.8
$ws.Range("J140").Value = 71764.22
$ws.Range("L140").Value = 71764.22
$ws.Range("N140").Value = -82124.22

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1245.3334
$ws.Range("I16").Value = 1368
$ws.Range("K16").Value = 1368
$ws.Range("M16").Value = -1081
$ws.Range("H45").Value = 7814431
$ws.Range("I45").Value = 2059.6365
$ws.Range("J45").Value = 25001648
$ws.Range("K45").Value = 2059.6365
$ws.Range("L45").Value = 25001648
$ws.Range("M45").Value = -1682.6365
$ws.Range("N45").Value = -25002402
$ws.Range("H74").Value = 24741.883
$ws.Range("I74").Value = 28019
$ws.Range("J74").Value = 4533
$ws.Range("K74").Value = 28019
$ws.Range("L74").Value = 4533
$ws.Range("M74").Value = -27145
$ws.Range("N74").Value = -6281
$ws.Range("H77").Value = 24741.883
$ws.Range("I77").Value = 28019
$ws.Range("J77").Value = 4533
$ws.Range("K77").Value = 140095
$ws.Range("L77").Value = 22665
$ws.Range("M77").Value = -135727
$ws.Range("N77").Value = -31401
$ws.Range("H122").Value = 2911.2144
$ws.Range("I122").Value = 2724.8
$ws.Range("J122").Value = 3377.25
$ws.Range("K122").Value = 8174.400000000001
$ws.Range("L122").Value = 10131.75
$ws.Range("M122").Value = -5724.400000000001
$ws.Range("N122").Value = -15031.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2440.889
$ws.Range("I94").Value = 2458.923
$ws.Range("J94").Value = 2394
$ws.Range("K94").Value = 2458.923
$ws.Range("L94").Value = 2394
$ws.Range("M94").Value = -2007.923
$ws.Range("N94").Value = -3296
$ws.Range("H105").Value = 44386.875
$ws.Range("I105").Value = 57488.332
$ws.Range("J105").Value = 5082.5
$ws.Range("K105").Value = 57488.332
$ws.Range("L105").Value = 5082.5
$ws.Range("M105").Value = -55741.332
$ws.Range("N105").Value = -8576.5
$ws.Range("H134").Value = 2624.5518
$ws.Range("I134").Value = 1291.0476
$ws.Range("J134").Value = 6125
$ws.Range("K134").Value = 3873.142800000001
$ws.Range("L134").Value = 18375
$ws.Range("M134").Value = -1338.142800000001
$ws.Range("N134").Value = -23445
$ws.Range("H140").Value = 45240.805
$ws.Range("I140").Value = 48845
$ws.Range("K140").Value = 48845
$ws.Range("M140").Value = -43665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2249.9473
$ws.Range("I31").Value = 1538.1724
$ws.Range("J31").Value = 2987.1428
$ws.Range("K31").Value = 1538.1724
$ws.Range("L31").Value = 2987.1428
$ws.Range("M31").Value = -1243.1724
$ws.Range("N31").Value = -3577.1428
$ws.Range("H34").Value = 2249.9473
$ws.Range("I34").Value = 1538.1724
$ws.Range("J34").Value = 2987.1428
$ws.Range("K34").Value = 1538.1724
$ws.Range("L34").Value = 2987.1428
$ws.Range("M34").Value = -1336.1724
$ws.Range("N34").Value = -3391.1428
$ws.Range("H94").Value = 707.6923
$ws.Range("J94").Value = 626.2222
$ws.Range("L94").Value = 626.2222
$ws.Range("N94").Value = -1528.2222
$ws.Range("H132").Value = 2379.75
$ws.Range("I132").Value = 2276.1304
$ws.Range("J132").Value = 2644.5557
$ws.Range("K132").Value = 6828.3912
$ws.Range("L132").Value = 7933.6671
$ws.Range("M132").Value = -4298.3912
$ws.Range("N132").Value = -12993.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 143.75
$ws.Range("I2").Value = 158.33333
$ws.Range("K2").Value = 949.9999799999999
$ws.Range("M2").Value = -836.9999799999999
$ws.Range("H4").Value = 6454661
$ws.Range("I4").Value = 7693153.5
$ws.Range("J4").Value = 14500.4
$ws.Range("K4").Value = 23079460.5
$ws.Range("L4").Value = 43501.2
$ws.Range("M4").Value = -23079348.5
$ws.Range("N4").Value = -43725.2
$ws.Range("H107").Value = 308
$ws.Range("I107").Value = 363.75
$ws.Range("K107").Value = 1091.25
$ws.Range("M107").Value = 828.75
$ws.Range("H113").Value = 126259.25
$ws.Range("J113").Value = 201354.8
$ws.Range("L113").Value = 604064.3999999999
$ws.Range("N113").Value = -608404.3999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 30000
$ws.Range("I34").Value = 30000
$ws.Range("K34").Value = 30000
$ws.Range("M34").Value = -29732
$ws.Range("H64").Value = 75000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
$ws.Range("H67").Value = 75000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
$ws.Range("H76").Value = 30000
$ws.Range("I76").Value = 30000
$ws.Range("K76").Value = 30000
$ws.Range("M76").Value = -29685
$ws.Range("H79").Value = 30000
$ws.Range("I79").Value = 30000
$ws.Range("K79").Value = 30000
$ws.Range("M79").Value = -28908
$ws.Range("H102").Value = 1782.6364
$ws.Range("I102").Value = 1460.9
$ws.Range("K102").Value = 1460.9
$ws.Range("M102").Value = 161.0999999999999
$ws.Range("H122").Value = 71298.35000000001
$ws.Range("I122").Value = 82767.14
$ws.Range("J122").Value = 4779.4
$ws.Range("K122").Value = 248301.42
$ws.Range("L122").Value = 14338.2
$ws.Range("M122").Value = -245851.42
$ws.Range("N122").Value = -19238.2
$ws.Range("H135").Value = 49664.777
$ws.Range("J135").Value = 49664.777
$ws.Range("L135").Value = 49664.777
$ws.Range("N135").Value = -59804.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22243.357
$ws.Range("I7").Value = 29340.7
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 29340.7
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -29228.7
$ws.Range("N7").Value = -4724
$ws.Range("H93").Value = 1738.7142
$ws.Range("I93").Value = 1809.2
$ws.Range("K93").Value = 1809.2
$ws.Range("M93").Value = -561.2
$ws.Range("H126").Value = 22243.357
$ws.Range("I126").Value = 29340.7
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 88022.10000000001
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -85552.10000000001
$ws.Range("N126").Value = -18440
$ws.Range("H132").Value = 1901.84
$ws.Range("I132").Value = 1716.0526
$ws.Range("J132").Value = 2490.1667
$ws.Range("K132").Value = 5148.1578
$ws.Range("L132").Value = 7470.500100000001
$ws.Range("M132").Value = -2618.1578
$ws.Range("N132").Value = -12530.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 8808.647000000001
$ws.Range("I107").Value = 10576.818
$ws.Range("J107").Value = 5567
$ws.Range("K107").Value = 31730.454
$ws.Range("L107").Value = 16701
$ws.Range("M107").Value = -29810.454
$ws.Range("N107").Value = -20541
$ws.Range("H113").Value = 1375.8182
$ws.Range("I113").Value = 1375.8182
$ws.Range("K113").Value = 4127.4546
$ws.Range("M113").Value = -1957.4546
$ws.Range("H122").Value = 3044
$ws.Range("I122").Value = 2448.75
$ws.Range("K122").Value = 7346.25
$ws.Range("M122").Value = -4896.25
$ws.Range("H132").Value = 4832568.5
$ws.Range("I132").Value = 1466.6666
$ws.Range("J132").Value = 7248119
$ws.Range("K132").Value = 4399.9998
$ws.Range("L132").Value = 21744357
$ws.Range("M132").Value = -1869.9998
$ws.Range("N132").Value = -21749417
$ws.Range("H136").Value = 1952.5319
$ws.Range("I136").Value = 1763.1951
$ws.Range("J136").Value = 3246.3333
$ws.Range("K136").Value = 5289.5853
$ws.Range("L136").Value = 9738.999899999999
$ws.Range("M136").Value = -2739.5853
$ws.Range("N136").Value = -14838.9999
